$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$A2 = @'
Pipeline(steps=[('scaler', StandardScaler()), ('selector', 'passthrough'),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(max_depth=1,
                                                                    max_features='sqrt',
                                                                    min_samples_leaf=3,
                                                                    min_samples_split=4,
                                                                    random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("A2").Value = $A2
$ws.Range("B2").Value = 0.6476190476190476

$C2 = @'
{'scaler': StandardScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 3, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 1, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': None}
'@
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = 0.5000000000000001

$E2 = @'
[1 0 0 1 0 0 1 1 0 1 0 0]
'@
$ws.Range("E2").Value = $E2

$ws.Range("H2").Value = 0.8626428571428572
$ws.Range("I2").Value = 0.02324063327676608
$ws.Range("J2").Value = 0.569047619047619
$ws.Range("K2").Value = 0.0691453109814083

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$A3 = @'
Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', 'passthrough'),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',
                                                                    max_depth=6,
                                                                    max_features='log2',
                                                                    min_samples_leaf=5,
                                                                    min_samples_split=4,
                                                                    random_state=42),
                                   random_state=42))])
'@
$ws.Range("A3").Value = $A3
$ws.Range("B3").Value = 0.6285714285714287

$C3 = @'
{'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 6, 'model__estimator__criterion': 'gini', 'model__estimator__class_weight': 'balanced'}
'@
$ws.Range("C3").Value = $C3
$ws.Range("D3").Value = 0.6153846153846153

$E3 = @'
[1 0 1 0 0 0 0 1 1 0 1 1]
'@
$ws.Range("E3").Value = $E3

$F3 = @'
[0 1 1 0 1 0 1 1 0 0 1 1]
'@
$ws.Range("F3").Value = $F3

$ws.Range("H3").Value = 0.8356190476190476
$ws.Range("I3").Value = 0.02432303225698967
$ws.Range("J3").Value = 0.5485714285714285
$ws.Range("K3").Value = 0.06199233572023411

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$A4 = @'
Pipeline(steps=[('scaler', StandardScaler()), ('selector', 'passthrough'),
                ('model',
                 BaggingClassifier(estimator=DecisionTreeClassifier(criterion='entropy',
                                                                    max_depth=5,
                                                                    max_features='log2',
                                                                    min_samples_leaf=3,
                                                                    random_state=42),
                                   random_state=42))])
'@
$ws.Range("A4").Value = $A4
$ws.Range("B4").Value = 0.5809523809523809

$C4 = @'
{'scaler': StandardScaler(), 'model__n_estimators': 10, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 3, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 5, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': None}
'@
$ws.Range("C4").Value = $C4
$ws.Range("D4").Value = 0.7777777777777777

$E4 = @'
[1 0 1 1 1 1 0 1 0 1 0 1]
'@
$ws.Range("E4").Value = $E4

$F4 = @'
[1 1 1 1 1 1 0 0 1 1 1 1]
'@
$ws.Range("F4").Value = $F4

$ws.Range("H4").Value = 0.8550476190476191
$ws.Range("I4").Value = 0.02805800153413196
$ws.Range("J4").Value = 0.5033333333333334
$ws.Range("K4").Value = 0.08857182129619702

$wb.Save()
